# update auto code generator
# Inserts a new "display name" header row (Korean labels) and a new
# "MinSpeed" float column into the code-gen template sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Drop the stray utility cell at K3 that lived outside the real
#    table (it only existed to pad the old dimension to K9).
# ------------------------------------------------------------------
$ws.Range("K3").Clear()

# ------------------------------------------------------------------
# 2. Insert a new row above the old "Id/Name/Nickname/Number" row so
#    we can put Korean display names there. This shifts rows 4-9 down
#    to 5-10 and auto-adjusts the ISNUMBER() formulas + merged cells.
# ------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Row 2 only keeps the "key" cell now; drop the leftover blank B2:D2.
$ws.Range("B2:D2").Clear()

# ------------------------------------------------------------------
# 3. Build the two brand-new cell styles we need by tweaking copies
#    of existing formatted cells (so we reuse fonts/fills instead of
#    re-declaring them), then stash them off in scratch cells we will
#    wipe again at the end.
# ------------------------------------------------------------------

# 3a. Yellow-highlighted header style (for the new Korean label row).
#     Based on the existing centered "no border" style used by A2.
$yellowHeader = $ws.Range("ZZ1")
$ws.Range("A2").Copy()
$yellowHeader.PasteSpecial(-4122)
$yellowHeader.Interior.Color = 65535

# 3b. Bold/gray header style with only left+right borders (for the new
#     "MinSpeed" column header), based on the existing full-border
#     bold header style used by A4 (old row4, now row5).
$sideBorderHeader = $ws.Range("ZZ2")
$ws.Range("A5").Copy()
$sideBorderHeader.PasteSpecial(-4122)
$sideBorderHeader.Borders.Item(8).LineStyle = -4142
$sideBorderHeader.Borders.Item(9).LineStyle = -4142

# ------------------------------------------------------------------
# 4. Fill in the new Korean label row (row 4).
# ------------------------------------------------------------------
$ws.Range("A4").Value = "아이디"
$ws.Range("B4").Value = "이름"
$ws.Range("C4").Value = "닉네임"
$ws.Range("D4").Value = "번호"

$yellowHeader.Copy()
$ws.Range("A4:E4").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 5. Fill in the new column E (design / float / MinSpeed / data).
#    Order matches the authored edit so shared-string ids line up.
# ------------------------------------------------------------------
$ws.Range("E1").Value = "design"
$ws.Range("E3").Value = "float"
$ws.Range("E4").Value = "스피드 범위"
$ws.Range("E5").Value = "MinSpeed"
$ws.Range("E6").Value = 50.05
$ws.Range("E7").Value = 50.05
$ws.Range("E8").Value = 50.05
$ws.Range("E9").Value = 50.05
$ws.Range("E10").Value = 50.05

# Match formatting of the rest of each row for the new column.
$ws.Range("A1").Copy()
$ws.Range("E1").PasteSpecial(-4122)

$sideBorderHeader.Copy()
$ws.Range("E5").PasteSpecial(-4122)

# ------------------------------------------------------------------
# 6. Clean up scratch cells used to build the new styles.
# ------------------------------------------------------------------
$yellowHeader.Clear()
$sideBorderHeader.Clear()

# ------------------------------------------------------------------
# 7. Restore selection to match the authored workbook.
# ------------------------------------------------------------------
$ws.Range("E7").Select()
